$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Widen column E to match the new content (closest attainable value given the
# host's pixel-quantized ColumnWidth rounding; target stored width is 34.7109375)
$ws.Columns.Item(5).ColumnWidth = 33.8

# New values for column E, rows 7-10 (with wrap text style matching s="2" on E7/E8)
# Shared-string insertion order must match the target: songsheader(21), songlink(22),
# songlikes-linktext(23), songlikes-action(24) -- so write E8 before E7.
$ws.Range("E8").Value = "songsheader: text : Songs`r`nlinkedsongs:  href=""/songs/5"""
$ws.Range("E7").Value = "songlink: linktext :Songs`r`nimagenSinatra: src=""/images/sinatra.jpg"""
$ws.Range("E9").Value = "songlikes: linktext : id=""like"", <p>This song has been liked 4 times</p>"
$ws.Range("E10").Value = "songlikes: action=""/songs/3/like"" ,id=""like"", <p>This song has been liked 4 times</p>"

$ws.Range("E7:E8").WrapText = $true

# Update the active selection to match the diff
$ws.Range("B8").Select()
